# Add season-record columns (Wins / Losses / Ties) to the stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): copy the existing header style (bold, bordered,
# centered) from AC1 onto the new AD1:AF1 header cells, then set labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-46: every player shares the team's overall season record.
$wins = 87
$losses = 75
$ties = 0

for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
